$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$formula = "=IF(ISBLANK(Tableau1[[#This Row],[Heure Début]]),`"`",Tableau1[[#This Row],[Heure fin]]-Tableau1[[#This Row],[Heure Début]])"

# Copy formatting from an already-filled row down onto the three blank rows
# so date/time/text styles match the rest of the table.
$ws.Range("E6:M6").Copy()
$ws.Range("E38:M38").PasteSpecial(-4122)
$ws.Range("E39:M39").PasteSpecial(-4122)
$ws.Range("E40:M40").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 38
$ws.Range("E38").Value = 44273
$ws.Range("F38").Value = 0.33333333333333331
$ws.Range("G38").Value = 0.35416666666666669
$ws.Range("H38").Formula = $formula
$ws.Range("I38").Value = "exercice"
$ws.Range("J38").Value = "touver la meilleur manière de tester un jeux"
$ws.Range("K38").Value = "CPNV"
$ws.Range("L38").Value = "Trouver la meilleur manière de tester un jeux tout en parlant en anglais à deux"

# Row 39
$ws.Range("E39").Value = 44273
$ws.Range("F39").Value = 0.35416666666666669
$ws.Range("G39").Value = 0.38541666666666669
$ws.Range("H39").Formula = $formula
$ws.Range("I39").Value = "théorie"
$ws.Range("J39").Value = "théorie sur les testes"
$ws.Range("K39").Value = "CPNV"
$ws.Range("L39").Value = "Théorie sur les différentes manière de tester un programme"
$ws.Range("M39").Value = "LVT"

# Row 40
$ws.Range("E40").Value = 44273
$ws.Range("F40").Value = 0.42708333333333331
$ws.Range("G40").Value = 0.44791666666666669
$ws.Range("H40").Formula = $formula
$ws.Range("I40").Value = "Documentation"
$ws.Range("J40").Value = "Scénario de teste"
$ws.Range("K40").Value = "CPNV"
$ws.Range("L40").Value = "faire un tableau excel pour les scénario de testes"

# Row heights (match the author's manual formatting in the saved file)
$ws.Range("E38:M38").RowHeight = 60.6
$ws.Range("E39:M39").RowHeight = 43.2
$ws.Range("E40:M40").RowHeight = 28.8

# Scroll position / selection left by Excel after the last edit
[void]$ws.Range("L41").Select()
$ws.Application.ActiveWindow.ScrollRow = 34
